{"js": "// Apply the Goldify copy-edit pass:\n//  - Shorten the H1 title (and the matching bold \"title\" run near the foot of\n//    the doc) to \"Play Goldify for Free\".\n//  - Rewrite the three \"What we like\" bullets and the first \"What we don't\n//    like\" bullet.\n//  - Rewrite the italic meta-description run.\nconst replacements = [\n  [\n    \"Play Goldify Slot for Free - Unique Symbol Gilding Feature\",\n    \"Play Goldify for Free\",\n  ],\n  [\"Immersive Ancient Greek theme\", \"Bonus round with up to 100 free spins\"],\n  [\"Lucrative bonus round with free spins\", \"Appealing Greek-themed graphics\"],\n  [\n    \"Graphics-rich and well-designed symbols\",\n    \"Potential for big wins with wild and scatter symbols\",\n  ],\n  [\n    \"Slightly below average payout range\",\n    \"Payout range slightly below average\",\n  ],\n  [\n    \"Experience the immersive Ancient Greek theme & lucrative bonus round with free spins. Play Goldify slot for free with unique symbol gilding feature.\",\n    \"Read our review of Goldify, a Greek-themed slot game with unique features. Play for free.\",\n  ],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Goldify copy-edit pass:\n#  - Shorten the H1 title (and the matching bold \"title\" run near the foot of\n#    the doc) to \"Play Goldify for Free\".\n#  - Rewrite the three \"What we like\" bullets and the first \"What we don't\n#    like\" bullet.\n#  - Rewrite the italic meta-description run.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"Play Goldify Slot for Free - Unique Symbol Gilding Feature\", \"Play Goldify for Free\"),\n    @(\"Immersive Ancient Greek theme\", \"Bonus round with up to 100 free spins\"),\n    @(\"Lucrative bonus round with free spins\", \"Appealing Greek-themed graphics\"),\n    @(\"Graphics-rich and well-designed symbols\", \"Potential for big wins with wild and scatter symbols\"),\n    @(\"Slightly below average payout range\", \"Payout range slightly below average\"),\n    @(\"Experience the immersive Ancient Greek theme & lucrative bonus round with free spins. Play Goldify slot for free with unique symbol gilding feature.\", \"Read our review of Goldify, a Greek-themed slot game with unique features. Play for free.\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
